$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4169.5386
$ws.Range("I86").Value = 1950
$ws.Range("J86").Value = 4573.091
$ws.Range("K86").Value = 1950
$ws.Range("L86").Value = 4573.091
$ws.Range("M86").Value = -827
$ws.Range("N86").Value = -6819.091

$ws.Range("H89").Value = 4169.5386
$ws.Range("I89").Value = 1950
$ws.Range("J89").Value = 4573.091
$ws.Range("K89").Value = 9750
$ws.Range("L89").Value = 22865.455
$ws.Range("M89").Value = -4134
$ws.Range("N89").Value = -34097.455

$ws.Range("H92").Value = 561.92
$ws.Range("I92").Value = 530.2143
$ws.Range("J92").Value = 602.2727
$ws.Range("K92").Value = 530.2143
$ws.Range("L92").Value = 602.2727
$ws.Range("M92").Value = 717.7857
$ws.Range("N92").Value = -3098.2727

$ws.Range("H100").Value = 1075
$ws.Range("I100").Value = 1166.6666
$ws.Range("J100").Value = 800
$ws.Range("K100").Value = 1166.6666
$ws.Range("L100").Value = 800
$ws.Range("M100").Value = -625.6666
$ws.Range("N100").Value = -1882

$ws.Range("H107").Value = 890
$ws.Range("I107").Value = 1378
$ws.Range("J107").Value = 585
$ws.Range("K107").Value = 1378
$ws.Range("L107").Value = 585
$ws.Range("M107").Value = 542
$ws.Range("N107").Value = -4425

$ws.Range("H112").Value = 2304.1707
$ws.Range("I112").Value = 967.75
$ws.Range("J112").Value = 2448.6487
$ws.Range("K112").Value = 2903.25
$ws.Range("L112").Value = 7345.946100000001
$ws.Range("M112").Value = -1795.25
$ws.Range("N112").Value = -9561.946100000001

$ws.Range("H116").Value = 1633.4546
$ws.Range("I116").Value = 1925
$ws.Range("J116").Value = 1466.8572
$ws.Range("K116").Value = 1925
$ws.Range("L116").Value = 1466.8572
$ws.Range("M116").Value = 1517
$ws.Range("N116").Value = -8350.8572

$ws.Range("H129").Value = 374.25
$ws.Range("I129").Value = 374.25
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 1122.75
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 3877.25
$ws.Range("N129").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8372.277
$ws.Range("I32").Value = 4629.19
$ws.Range("J32").Value = 35254.453
$ws.Range("K32").Value = 4629.19
$ws.Range("L32").Value = 35254.453
$ws.Range("M32").Value = -4342.19
$ws.Range("N32").Value = -35828.453

$ws.Range("H61").Value = 3473790.5
$ws.Range("I61").Value = 3969791.5
$ws.Range("J61").Value = 1782
$ws.Range("K61").Value = 3969791.5
$ws.Range("L61").Value = 1782
$ws.Range("M61").Value = -3969579.5
$ws.Range("N61").Value = -2206

$ws.Range("H97").Value = 977.5
$ws.Range("I97").Value = 955
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 955
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -459
$ws.Range("N97").Value = -1992

$ws.Range("H122").Value = 1348.9656
$ws.Range("I122").Value = 1160.4375
$ws.Range("J122").Value = 1581
$ws.Range("K122").Value = 3481.3125
$ws.Range("L122").Value = 4743
$ws.Range("M122").Value = -1031.3125
$ws.Range("N122").Value = -9643

$ws.Range("H132").Value = 936531.9399999999
$ws.Range("I132").Value = 1382208
$ws.Range("J132").Value = 79462.46000000001
$ws.Range("K132").Value = 4146624
$ws.Range("L132").Value = 238387.38
$ws.Range("M132").Value = -4144094
$ws.Range("N132").Value = -243447.38

$ws.Range("H136").Value = 3473790.5
$ws.Range("I136").Value = 3969791.5
$ws.Range("J136").Value = 1782
$ws.Range("K136").Value = 11909374.5
$ws.Range("L136").Value = 5346
$ws.Range("M136").Value = -11906824.5
$ws.Range("N136").Value = -10446

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 957.8
$ws.Range("I94").Value = 866.7
$ws.Range("J94").Value = 1140
$ws.Range("K94").Value = 866.7
$ws.Range("L94").Value = 1140
$ws.Range("M94").Value = -415.7
$ws.Range("N94").Value = -2042

$ws.Range("H134").Value = 18274072
$ws.Range("I134").Value = 18274072
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 54822216
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -54819681

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 416.66666
$ws.Range("I105").Value = 375
$ws.Range("J105").Value = 500
$ws.Range("K105").Value = 375
$ws.Range("L105").Value = 500
$ws.Range("M105").Value = 1372
$ws.Range("N105").Value = -3994

$ws.Range("H122").Value = 1273.9333
$ws.Range("I122").Value = 1060
$ws.Range("J122").Value = 1380.9
$ws.Range("K122").Value = 3180
$ws.Range("L122").Value = 4142.700000000001
$ws.Range("M122").Value = -730
$ws.Range("N122").Value = -9042.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 2576.4707
$ws.Range("I56").Value = 2576.4707
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 2576.4707
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -2046.4707

$ws.Range("H131").Value = 1945.6495
$ws.Range("I131").Value = 3647.3125
$ws.Range("J131").Value = 1609.5186
$ws.Range("K131").Value = 10941.9375
$ws.Range("L131").Value = 4828.5558
$ws.Range("M131").Value = -5901.9375
$ws.Range("N131").Value = -14908.5558

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2085.1
$ws.Range("I97").Value = 2034.6666
$ws.Range("J97").Value = 2236.4
$ws.Range("K97").Value = 2034.6666
$ws.Range("L97").Value = 2236.4
$ws.Range("M97").Value = -1538.6666
$ws.Range("N97").Value = -3228.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2414.7856
$ws.Range("I93").Value = 2379.9
$ws.Range("J93").Value = 2502
$ws.Range("K93").Value = 2379.9
$ws.Range("L93").Value = 2502
$ws.Range("M93").Value = -1131.9
$ws.Range("N93").Value = -4998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 16796.75
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 16796.75
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 16796.75
$ws.Range("N46").Value = -17258.75

$ws.Range("H107").Value = 282.8
$ws.Range("I107").Value = 270.18182
$ws.Range("J107").Value = 317.5
$ws.Range("K107").Value = 810.54546
$ws.Range("L107").Value = 952.5
$ws.Range("M107").Value = 1109.45454
$ws.Range("N107").Value = -4792.5

$ws.Range("H113").Value = 655.1905
$ws.Range("I113").Value = 611.5714
$ws.Range("J113").Value = 677
$ws.Range("K113").Value = 1834.7142
$ws.Range("L113").Value = 2031
$ws.Range("M113").Value = 335.2857999999999
$ws.Range("N113").Value = -6371

$ws.Range("H122").Value = 1183.5
$ws.Range("I122").Value = 1237.5714
$ws.Range("J122").Value = 1107.8
$ws.Range("K122").Value = 3712.7142
$ws.Range("L122").Value = 3323.4
$ws.Range("M122").Value = -1262.7142
$ws.Range("N122").Value = -8223.4

$ws.Range("H134").Value = 16796.75
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 16796.75
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 50390.25
$ws.Range("N134").Value = -55460.25
